# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 1657
$ws1.Range("F3").Value  = 9074
$ws1.Range("F4").Value  = 110
$ws1.Range("F5").Value  = 505
$ws1.Range("F7").Value  = 1232
$ws1.Range("F10").Value = 91
$ws1.Range("F11").Value = 5878
$ws1.Range("F13").Value = 384
$ws1.Range("F15").Value = 4384
$ws1.Range("F16").Value = 10
$ws1.Range("F17").Value = 162
$ws1.Range("F18").Value = 1147
$ws1.Range("F19").Value = 22
$ws1.Range("F24").Value = 2722

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 1657
$ws4.Range("F3").Value  = 9074
$ws4.Range("F4").Value  = 110
$ws4.Range("F6").Value  = 505
$ws4.Range("F8").Value  = 1233
$ws4.Range("F11").Value = 91
$ws4.Range("F12").Value = 5878
$ws4.Range("F14").Value = 384
$ws4.Range("F16").Value = 4384
$ws4.Range("F17").Value = 10
$ws4.Range("F18").Value = 162
$ws4.Range("F19").Value = 1147
$ws4.Range("F20").Value = 22
$ws4.Range("F25").Value = 2722
